$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.3611711538166844
$ws.Cells.Item(2, 4).Value = 0.05293986698581676
$ws.Cells.Item(2, 5).Value = 0.1341110833114421
$ws.Cells.Item(2, 6).Value = 3.414142300303382
$ws.Cells.Item(2, 7).Value = 2.674864574869062
$ws.Cells.Item(2, 8).Value = 2.144191772653727
$ws.Cells.Item(2, 10).Value = 0.2521415955952264
$ws.Cells.Item(2, 11).Value = 2.278151364513633
$ws.Cells.Item(2, 12).Value = 0.1073198676380862
$ws.Cells.Item(2, 13).Value = 0.6818063226443911
$ws.Cells.Item(3, 3).Value = 0.3593735820300168
$ws.Cells.Item(3, 4).Value = 0.05203751109750243
$ws.Cells.Item(3, 5).Value = 0.1343789104677509
$ws.Cells.Item(3, 6).Value = 3.428192868062069
$ws.Cells.Item(3, 7).Value = 2.686204409487985
$ws.Cells.Item(3, 8).Value = 2.157466703631442
$ws.Cells.Item(3, 10).Value = 0.2538375658741359
$ws.Cells.Item(3, 11).Value = 2.156780224461102
$ws.Cells.Item(3, 12).Value = 0.1076747805200879
$ws.Cells.Item(3, 13).Value = 0.6601000140164786
$ws.Cells.Item(4, 3).Value = 0.3584268674124189
$ws.Cells.Item(4, 4).Value = 0.05149047033647491
$ws.Cells.Item(4, 5).Value = 0.1345842861272502
$ws.Cells.Item(4, 6).Value = 3.438640194198314
$ws.Cells.Item(4, 7).Value = 2.694775964570951
$ws.Cells.Item(4, 8).Value = 2.166647454514745
$ws.Cells.Item(4, 10).Value = 0.2549698984456157
$ws.Cells.Item(4, 11).Value = 2.08308952106637
$ws.Cells.Item(4, 12).Value = 0.1079077949471143
$ws.Cells.Item(4, 13).Value = 0.6470419749549166
$ws.Cells.Item(5, 3).Value = 0.358080582085563
$ws.Cells.Item(5, 4).Value = 0.05126933537287215
$ws.Cells.Item(5, 5).Value = 0.1346782885769464
$ws.Cells.Item(5, 6).Value = 3.443354698317805
$ws.Cells.Item(5, 7).Value = 2.698672673456571
$ws.Cells.Item(5, 8).Value = 2.170647445334893
$ws.Cells.Item(5, 10).Value = 0.2554542238731656
$ws.Cells.Item(5, 11).Value = 2.053270076689074
$ws.Cells.Item(5, 12).Value = 0.1080065565083403
$ws.Cells.Item(5, 13).Value = 0.6417887674637228
$ws.Cells.Item(6, 3).Value = 0.3580254687772992
$ws.Cells.Item(6, 4).Value = 0.05123272488373232
$ws.Cells.Item(6, 5).Value = 0.1346945208662422
$ws.Cells.Item(6, 6).Value = 3.444165131156311
$ws.Cells.Item(6, 7).Value = 2.699344079167076
$ws.Cells.Item(6, 8).Value = 2.171327263040752
$ws.Cells.Item(6, 10).Value = 0.2555360285247623
$ws.Cells.Item(6, 11).Value = 2.048331300973871
$ws.Cells.Item(6, 12).Value = 0.1080231859798655
$ws.Cells.Item(6, 13).Value = 0.6409205931802475
$ws.Cells.Item(7, 3).Value = 0.3584220372774212
$ws.Cells.Item(7, 4).Value = 0.05148748075512088
$ws.Cells.Item(7, 5).Value = 0.1345855121068436
$ws.Cells.Item(7, 6).Value = 3.438701925512092
$ws.Cells.Item(7, 7).Value = 2.694826883371888
$ws.Cells.Item(7, 8).Value = 2.166700352377362
$ws.Cells.Item(7, 10).Value = 0.2549763375342486
$ws.Cells.Item(7, 11).Value = 2.082686513531399
$ws.Cells.Item(7, 12).Value = 0.1079091114551414
$ws.Cells.Item(7, 13).Value = 0.6469708525286606
$ws.Cells.Item(8, 3).Value = 0.3605187811576371
$ws.Cells.Item(8, 4).Value = 0.05262729499241914
$ws.Cells.Item(8, 5).Value = 0.1341949461646834
$ws.Cells.Item(8, 6).Value = 3.418608797240438
$ws.Cells.Item(8, 7).Value = 2.678440129890589
$ws.Cells.Item(8, 8).Value = 2.148555084728883
$ws.Cells.Item(8, 10).Value = 0.2527074838749002
$ws.Cells.Item(8, 11).Value = 2.236130521950429
$ws.Cells.Item(8, 12).Value = 0.1074391144261613
$ws.Cells.Item(8, 13).Value = 0.6742661316943099
$ws.Cells.Item(9, 3).Value = 0.3658756107919032
$ws.Cells.Item(9, 4).Value = 0.05491709295650082
$ws.Cells.Item(9, 5).Value = 0.1337530661593878
$ws.Cells.Item(9, 6).Value = 3.393677389327308
$ws.Cells.Item(9, 7).Value = 2.659110739866094
$ws.Cells.Item(9, 8).Value = 2.121154504863739
$ws.Cells.Item(9, 10).Value = 0.2489799851030412
$ws.Cells.Item(9, 11).Value = 2.543608239427272
$ws.Cells.Item(9, 12).Value = 0.1066367822052445
$ws.Cells.Item(9, 13).Value = 0.7299260046790508
$ws.Cells.Item(10, 3).Value = 0.3705704932402227
$ws.Cells.Item(10, 4).Value = 0.05663160128567313
$ws.Cells.Item(10, 5).Value = 0.1336250723430084
$ws.Cells.Item(10, 6).Value = 3.384225512252286
$ws.Cells.Item(10, 7).Value = 2.652773442196349
$ws.Cells.Item(10, 8).Value = 2.106027267355074
$ws.Cells.Item(10, 10).Value = 0.2466809685736244
$ws.Cells.Item(10, 11).Value = 2.773515241644247
$ws.Cells.Item(10, 12).Value = 0.1061194358694246
$ws.Cells.Item(10, 13).Value = 0.7721173340004555
$ws.Cells.Item(11, 3).Value = 0.3728712959331233
$ws.Cells.Item(11, 4).Value = 0.05741835431084752
$ws.Cells.Item(11, 5).Value = 0.133609372656128
$ws.Cells.Item(11, 6).Value = 3.381860287215218
$ws.Cells.Item(11, 7).Value = 2.651610768407608
$ws.Cells.Item(11, 8).Value = 2.100235641360911
$ws.Cells.Item(11, 10).Value = 0.2457304623973968
$ws.Cells.Item(11, 11).Value = 2.87897648366453
$ws.Cells.Item(11, 12).Value = 0.1058996125248131
$ws.Cells.Item(11, 13).Value = 0.7915929778486657
$ws.Cells.Item(12, 3).Value = 0.3737662815246665
$ws.Cells.Item(12, 4).Value = 0.05771723616989277
$ws.Cells.Item(12, 5).Value = 0.1336095271117621
$ws.Cells.Item(12, 6).Value = 3.381243559780046
$ws.Cells.Item(12, 7).Value = 2.651418854739632
$ws.Cells.Item(12, 8).Value = 2.098199508124338
$ws.Cells.Item(12, 10).Value = 0.2453842340356829
$ws.Cells.Item(12, 11).Value = 2.919037381691965
$ws.Cells.Item(12, 12).Value = 0.1058185927424837
$ws.Cells.Item(12, 13).Value = 0.7990084227240004
$ws.Cells.Item(13, 3).Value = 0.3735724755884462
$ws.Cells.Item(13, 4).Value = 0.05765282449607412
$ws.Cells.Item(13, 5).Value = 0.1336092228104171
$ws.Cells.Item(13, 6).Value = 3.381363965374064
$ws.Cells.Item(13, 7).Value = 2.65144912459624
$ws.Cells.Item(13, 8).Value = 2.098631036863196
$ws.Cells.Item(13, 10).Value = 0.2454581908562687
$ws.Cells.Item(13, 11).Value = 2.910404003837186
$ws.Cells.Item(13, 12).Value = 0.1058359431001827
$ws.Cells.Item(13, 13).Value = 0.7974095783903152
$ws.Cells.Item(14, 3).Value = 0.3729444516092428
$ws.Cells.Item(14, 4).Value = 0.05744292448747501
$ws.Cells.Item(14, 5).Value = 0.1336092632062105
$ws.Cells.Item(14, 6).Value = 3.381803954111305
$ws.Cells.Item(14, 7).Value = 2.651589996165228
$ws.Cells.Item(14, 8).Value = 2.100064978532487
$ws.Cells.Item(14, 10).Value = 0.245701703260039
$ws.Cells.Item(14, 11).Value = 2.882269817200267
$ws.Cells.Item(14, 12).Value = 0.1058929024884714
$ws.Cells.Item(14, 13).Value = 0.7922022415527152
$ws.Cells.Item(15, 3).Value = 0.3725628575726603
$ws.Cells.Item(15, 4).Value = 0.05731447840202009
$ws.Cells.Item(15, 5).Value = 0.1336100818404802
$ws.Cells.Item(15, 6).Value = 3.382109807647694
$ws.Cells.Item(15, 7).Value = 2.651708658506806
$ws.Cells.Item(15, 8).Value = 2.100963768502908
$ws.Cells.Item(15, 10).Value = 0.2458526467812554
$ws.Cells.Item(15, 11).Value = 2.865053071071429
$ws.Cells.Item(15, 12).Value = 0.1059280809266658
$ws.Cells.Item(15, 13).Value = 0.7890178583450762
$ws.Cells.Item(16, 3).Value = 0.370423450719997
$ws.Cells.Item(16, 4).Value = 0.05658032015492864
$ws.Cells.Item(16, 5).Value = 0.1336269527416327
$ws.Cells.Item(16, 6).Value = 3.384419076429651
$ws.Cells.Item(16, 7).Value = 2.652884126037407
$ws.Cells.Item(16, 8).Value = 2.106427720616438
$ws.Cells.Item(16, 10).Value = 0.2467450045924267
$ws.Cells.Item(16, 11).Value = 2.766640641683637
$ws.Cells.Item(16, 12).Value = 0.1061341133560356
$ws.Cells.Item(16, 13).Value = 0.7708502245196769
$ws.Cells.Item(17, 3).Value = 0.3691532620760114
$ws.Cells.Item(17, 4).Value = 0.0561316657019475
$ws.Cells.Item(17, 5).Value = 0.1336481831379217
$ws.Cells.Item(17, 6).Value = 3.386331740603737
$ws.Cells.Item(17, 7).Value = 2.654046532597761
$ws.Cells.Item(17, 8).Value = 2.110059031729321
$ws.Cells.Item(17, 10).Value = 0.2473168526258007
$ws.Cells.Item(17, 11).Value = 2.706491384466801
$ws.Cells.Item(17, 12).Value = 0.106264476060419
$ws.Cells.Item(17, 13).Value = 0.7597772075584572
$ws.Cells.Item(18, 3).Value = 0.3684382219010445
$ws.Cells.Item(18, 4).Value = 0.05587425436667814
$ws.Cells.Item(18, 5).Value = 0.1336643978985013
$ws.Cells.Item(18, 6).Value = 3.387613891405962
$ws.Cells.Item(18, 7).Value = 2.654876970636849
$ws.Cells.Item(18, 8).Value = 2.112250226853803
$ws.Cells.Item(18, 10).Value = 0.247654737211807
$ws.Cells.Item(18, 11).Value = 2.671977568091677
$ws.Cells.Item(18, 12).Value = 0.1063409186118136
$ws.Cells.Item(18, 13).Value = 0.753434911610853
$ws.Cells.Item(19, 3).Value = 0.3681987908968836
$ws.Cells.Item(19, 4).Value = 0.05578721047211843
$ws.Cells.Item(19, 5).Value = 0.133670576036117
$ws.Cells.Item(19, 6).Value = 3.388079247722459
$ws.Cells.Item(19, 7).Value = 2.655185910514632
$ws.Cells.Item(19, 8).Value = 2.113009733322883
$ws.Cells.Item(19, 10).Value = 0.2477706802157869
$ws.Cells.Item(19, 11).Value = 2.660305970018385
$ws.Cells.Item(19, 12).Value = 0.106367052020687
$ws.Cells.Item(19, 13).Value = 0.7512920956574334
$ws.Cells.Item(20, 3).Value = 0.3692868677480732
$ws.Cells.Item(20, 4).Value = 0.05617935938108332
$ws.Cells.Item(20, 5).Value = 0.1336455088822603
$ws.Cells.Item(20, 6).Value = 3.386109288706052
$ws.Cells.Item(20, 7).Value = 2.653906034107848
$ws.Cells.Item(20, 8).Value = 2.109661855565975
$ws.Cells.Item(20, 10).Value = 0.2472550497894872
$ws.Cells.Item(20, 11).Value = 2.712885849187899
$ws.Cells.Item(20, 12).Value = 0.1062504475482413
$ws.Cells.Item(20, 13).Value = 0.7609531964733947
$ws.Cells.Item(21, 3).Value = 0.3731282738031041
$ws.Cells.Item(21, 4).Value = 0.05750455146615252
$ws.Cells.Item(21, 5).Value = 0.1336090859222452
$ws.Cells.Item(21, 6).Value = 3.381667142558911
$ws.Cells.Item(21, 7).Value = 2.651541870166028
$ws.Cells.Item(21, 8).Value = 2.099639530840165
$ws.Cells.Item(21, 10).Value = 0.2456298057420483
$ws.Cells.Item(21, 11).Value = 2.890530123647807
$ws.Cells.Item(21, 12).Value = 0.1058761118889517
$ws.Cells.Item(21, 13).Value = 0.7937306665679102
$ws.Cells.Item(22, 3).Value = 0.3757771154136833
$ws.Cells.Item(22, 4).Value = 0.05837619996906795
$ws.Cells.Item(22, 5).Value = 0.1336208274592217
$ws.Cells.Item(22, 6).Value = 3.380390128263926
$ws.Cells.Item(22, 7).Value = 2.651444887991659
$ws.Cells.Item(22, 8).Value = 2.094004800789747
$ws.Cells.Item(22, 10).Value = 0.2446475138549715
$ws.Cells.Item(22, 11).Value = 3.007359232381248
$ws.Cells.Item(22, 12).Value = 0.1056444123758311
$ws.Cells.Item(22, 13).Value = 0.8153882752044694
$ws.Cells.Item(23, 3).Value = 0.3743507329854481
$ws.Cells.Item(23, 4).Value = 0.05791048401743382
$ws.Cells.Item(23, 5).Value = 0.1336113133695314
$ws.Cells.Item(23, 6).Value = 3.380922647071799
$ws.Cells.Item(23, 7).Value = 2.65136381132919
$ws.Cells.Item(23, 8).Value = 2.096928294871304
$ws.Cells.Item(23, 10).Value = 0.2451644707739469
$ws.Cells.Item(23, 11).Value = 2.944938974209208
$ws.Cells.Item(23, 12).Value = 0.1057668928248674
$ws.Cells.Item(23, 13).Value = 0.8038077037232654
$ws.Cells.Item(24, 3).Value = 0.3692264172524773
$ws.Cells.Item(24, 4).Value = 0.05615779542877419
$ws.Cells.Item(24, 5).Value = 0.1336467054243684
$ws.Cells.Item(24, 6).Value = 3.38620929068415
$ws.Cells.Item(24, 7).Value = 2.653969048451756
$ws.Cells.Item(24, 8).Value = 2.109841096515993
$ws.Cells.Item(24, 10).Value = 0.247282962439261
$ws.Cells.Item(24, 11).Value = 2.70999470357367
$ws.Cells.Item(24, 12).Value = 0.1062567851801868
$ws.Cells.Item(24, 13).Value = 0.7604214579580031
$ws.Cells.Item(25, 3).Value = 0.3642931358397448
$ws.Cells.Item(25, 4).Value = 0.05429191065432803
$ws.Cells.Item(25, 5).Value = 0.1338380093858653
$ws.Cells.Item(25, 6).Value = 3.398868315159262
$ws.Cells.Item(25, 7).Value = 2.662963042027144
$ws.Cells.Item(25, 8).Value = 2.127689497642706
$ws.Cells.Item(25, 10).Value = 0.2499111502476268
$ws.Cells.Item(25, 11).Value = 2.459724175505414
$ws.Cells.Item(25, 12).Value = 0.1068411235360571
$ws.Cells.Item(25, 13).Value = 0.7146404101159263
